# Auto-generated Excel COM-interop script to apply meteocat data update
# Commit: Update automàtic: dades i banners [2026-02-09 21:50]
#
# Refreshes the per-station extraction timestamps (col E) together with the
# day's cumulative rainfall, humidity, snow depth, pressure and temperature
# readings (cols G/H/I/J/O) that changed between the 21:18-21:20 run and the
# 21:48-21:50 run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-09 21:48:56'
$ws.Range("I2").Value = '3.3 mm'
$ws.Range("E3").Value = '2026-02-09 21:48:59'
$ws.Range("G3").Value = '170 cm'
$ws.Range("I3").Value = '3.5 mm'
$ws.Range("E4").Value = '2026-02-09 21:49:01'
$ws.Range("E5").Value = '2026-02-09 21:49:04'
$ws.Range("G5").Value = '120 cm'
$ws.Range("E6").Value = '2026-02-09 21:49:07'
$ws.Range("E7").Value = '2026-02-09 21:49:09'
$ws.Range("I7").Value = '0.1 mm'
$ws.Range("J7").Value = '1007.6 hPa'
$ws.Range("E8").Value = '2026-02-09 21:49:12'
$ws.Range("I8").Value = '0.2 mm'
$ws.Range("J8").Value = '1007.5 hPa'
$ws.Range("E9").Value = '2026-02-09 21:49:15'
$ws.Range("E10").Value = '2026-02-09 21:49:17'
$ws.Range("O10").Value = '8.0 °C'
$ws.Range("E11").Value = '2026-02-09 21:49:20'
$ws.Range("H11").Value = '''83%'
$ws.Range("E12").Value = '2026-02-09 21:49:22'
$ws.Range("E13").Value = '2026-02-09 21:49:24'
$ws.Range("I13").Value = '0.2 mm'
$ws.Range("E14").Value = '2026-02-09 21:49:27'
$ws.Range("I14").Value = '0.2 mm'
$ws.Range("E15").Value = '2026-02-09 21:49:29'
$ws.Range("O15").Value = '8.1 °C'
$ws.Range("E16").Value = '2026-02-09 21:49:32'
$ws.Range("H16").Value = '''74%'
$ws.Range("I16").Value = '2.6 mm'
$ws.Range("E17").Value = '2026-02-09 21:49:35'
$ws.Range("E18").Value = '2026-02-09 21:49:37'
$ws.Range("E19").Value = '2026-02-09 21:49:40'
$ws.Range("E20").Value = '2026-02-09 21:49:43'
$ws.Range("I20").Value = '0.8 mm'
$ws.Range("O20").Value = '-4.2 °C'
$ws.Range("E21").Value = '2026-02-09 21:49:45'
$ws.Range("I21").Value = '0.6 mm'
$ws.Range("E22").Value = '2026-02-09 21:49:48'
$ws.Range("G22").Value = '120 cm'
$ws.Range("E23").Value = '2026-02-09 21:49:50'
$ws.Range("H23").Value = '''88%'
$ws.Range("I23").Value = '2.4 mm'
$ws.Range("E24").Value = '2026-02-09 21:49:53'
$ws.Range("I24").Value = '2.0 mm'
$ws.Range("J24").Value = '1008.5 hPa'
$ws.Range("E25").Value = '2026-02-09 21:49:56'
$ws.Range("H25").Value = '''75%'
$ws.Range("I25").Value = '0.7 mm'
$ws.Range("E26").Value = '2026-02-09 21:49:58'
$ws.Range("E27").Value = '2026-02-09 21:50:01'
$ws.Range("I27").Value = '1.0 mm'
$ws.Range("E28").Value = '2026-02-09 21:50:03'
$ws.Range("E29").Value = '2026-02-09 21:50:06'
$ws.Range("E30").Value = '2026-02-09 21:50:09'
$ws.Range("E31").Value = '2026-02-09 21:50:11'
$ws.Range("O31").Value = '9.8 °C'
$ws.Range("E32").Value = '2026-02-09 21:50:14'
$ws.Range("I32").Value = '1.8 mm'
$ws.Range("E33").Value = '2026-02-09 21:50:17'
$ws.Range("H33").Value = '''78%'
$ws.Range("I33").Value = '0.3 mm'
$ws.Range("E34").Value = '2026-02-09 21:50:19'
$ws.Range("H34").Value = '''76%'
$ws.Range("I34").Value = '0.1 mm'
$ws.Range("E35").Value = '2026-02-09 21:50:22'
$ws.Range("J35").Value = '1008.9 hPa'
$ws.Range("O35").Value = '5.5 °C'
$ws.Range("E36").Value = '2026-02-09 21:50:24'
$ws.Range("E37").Value = '2026-02-09 21:50:27'
$ws.Range("E38").Value = '2026-02-09 21:50:29'
$ws.Range("E39").Value = '2026-02-09 21:50:32'
$ws.Range("E40").Value = '2026-02-09 21:50:34'
$ws.Range("I40").Value = '1.1 mm'
$ws.Range("J40").Value = '1008.5 hPa'
$ws.Range("E41").Value = '2026-02-09 21:50:37'
$ws.Range("I41").Value = '0.7 mm'
$ws.Range("E42").Value = '2026-02-09 21:50:39'
$ws.Range("H42").Value = '''86%'
$ws.Range("E43").Value = '2026-02-09 21:50:42'
$ws.Range("E44").Value = '2026-02-09 21:50:45'
$ws.Range("I44").Value = '2.0 mm'
$ws.Range("O44").Value = '-3.8 °C'
$ws.Range("E45").Value = '2026-02-09 21:50:47'
$ws.Range("I45").Value = '1.3 mm'
$ws.Range("E46").Value = '2026-02-09 21:50:50'
$ws.Range("H46").Value = '''75%'
$ws.Range("J46").Value = '1008.7 hPa'
